$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1.00939432679166
$ws.Range("D2").Value = 0.07673537214160923
$ws.Range("C3").Value = 0.9749754488946234
$ws.Range("D3").Value = -0.05053452184021887
$ws.Range("F3").Value = 0.1028
$ws.Range("C4").Value = 0.9622904955496359
$ws.Range("D4").Value = -0.05520162107431915
$ws.Range("F4").Value = 0.0985
$ws.Range("C5").Value = 0.9531193669295966
$ws.Range("D5").Value = -0.09191377460253754
$ws.Range("F5").Value = 0.0963
$ws.Range("C6").Value = 0.9520514358513943
$ws.Range("D6").Value = -0.09682861819070865
$ws.Range("E6").Value = 0.9523
$ws.Range("F6").Value = 0.1038
$ws.Range("C7").Value = 0.9725155429972897
$ws.Range("D7").Value = -0.05837308947739608
$ws.Range("C8").Value = 0.9716480589416294
$ws.Range("D8").Value = -0.05991161094452264
$ws.Range("E8").Value = 0.9717
$ws.Range("C9").Value = 0.9673919320071813
$ws.Range("D9").Value = -0.07434942497228972
$ws.Range("F9").Value = 0.101
$ws.Range("C10").Value = 0.9643529542336976
$ws.Range("D10").Value = -0.08393696208884831
$ws.Range("E10").Value = 0.9643
$ws.Range("F10").Value = 0.1005
$ws.Range("C11").Value = 0.9658752548592051
$ws.Range("D11").Value = -0.07915722682329414
$ws.Range("E11").Value = 0.9659
$ws.Range("F11").Value = 0.0989
$ws.Range("C12").Value = 0.9575286848742928
$ws.Range("D12").Value = -0.06826733089280886
$ws.Range("F12").Value = 0.1
$ws.Range("C13").Value = 0.9554739975739441
$ws.Range("D13").Value = -0.07807269246896029
$ws.Range("E13").Value = 0.9550999999999999
$ws.Range("F13").Value = 0.09959999999999999
$ws.Range("C14").Value = 0.953759865174989
$ws.Range("D14").Value = -0.08620578960271544
$ws.Range("E14").Value = 0.9535
$ws.Range("F14").Value = 0.1017
$ws.Range("C15").Value = 0.9515221951387307
$ws.Range("D15").Value = -0.1075007037841103
$ws.Range("E15").Value = 0.9519
$ws.Range("F15").Value = 0.1055
$ws.Range("C16").Value = 0.9499201049895872
$ws.Range("D16").Value = -0.1230862792453686
$ws.Range("E16").Value = 0.9503
$ws.Range("F16").Value = 0.1107
